$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep only the hyperlink that lives on A2; drop the other three
# (they pointed at emails that are being removed from the sheet).
$keepAddr = '$A$2'
$changed = $true
while ($changed) {
  $changed = $false
  foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -ne $keepAddr) {
      $h.Delete()
      $changed = $true
      break
    }
  }
}

# Point the surviving hyperlink at the new recipient address
foreach ($h in $ws.Hyperlinks) {
  $h.Address = "mailto:k0vbasyuk.dim0n@gmail.com"
}

# Header row stays "email"; A2 becomes the single recipient email
$ws.Range("A1").Value = "email"
$ws.Range("A2").Value = "k0vbasyuk.dim0n@gmail.com"

# Rows 3-5 lose their old email values, but keep the hyperlink-style formatting
$ws.Range("A3:A5").ClearContents()

# Row 6 is a new blank row with the same (hyperlink) style as the rows above it
$ws.Range("A6").Style = "Гіперпосилання"

$ws.Range("A3").Select()
